# Apply the "ROM Map" update: double the size of three ROM chips
# (8192 -> 16384), fill in the bit-range / bit-pattern labels for rows
# 22-29, fix the mislabeled header, and leave the selection on K10 as the
# author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ROM Map")

# Rows 22-24: the actual ROM size doubled from 8192 to 16384 bytes.
$ws.Range("E22").Value() = 16384
$ws.Range("E23").Value() = 16384
$ws.Range("E24").Value() = 16384

# Fill in the bit-range ("L") column for rows 22-29 (leading apostrophe
# keeps these as literal text, same as the existing quote-prefixed cells
# above them, instead of being reinterpreted as a time value like 18:14).
$ws.Range("L22").Value() = "'18:14"
$ws.Range("L23").Value() = "'18:14"
$ws.Range("L24").Value() = "'18:14"
$ws.Range("L25").Value() = "'18:13"
$ws.Range("L26").Value() = "'18:13"
$ws.Range("L27").Value() = "'18:13"
$ws.Range("L28").Value() = "'18:12"

# Fill in the bit-pattern ("M") column for rows 22-24.
$ws.Range("M22").Value() = "01110"
$ws.Range("M23").Value() = "01111"
$ws.Range("M24").Value() = "10000"

# G1 was incorrectly labelled "d'Start" (same as F1); it should read "d'End".
$ws.Range("G1").Value() = "d'End"

# Row 29's bit range, then the remaining bit patterns for rows 25-29.
$ws.Range("L29").Value() = "'18:9"

$ws.Range("M25").Value() = "100010"
$ws.Range("M26").Value() = "100010"
$ws.Range("M27").Value() = "100011"
$ws.Range("M28").Value() = "1001000"
$ws.Range("M29").Value() = "1001001000"

# The conditional-formatting rules on J3:J29 and K3:K29 got renumbered by
# Excel during the edit session.
$ws.Range("J3:J29").FormatConditions.Item(1).Priority = 118
$ws.Range("K3:K29").FormatConditions.Item(1).Priority = 120

# Restore the active selection to K10, matching the author's saved view.
$ws.Activate()
$ws.Range("K10").Select()
